# Insert a new "Industry" column at column C, shifting Mutual Fund, Status,
# Jan_2026, Dec_2025, Oct_2025, MoM, QoQ one column to the right (D..J).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()

# Header
$ws.Range("C1").Value = "Industry"

# Row -> Industry value (matches the row order in the sheet)
$industries = @{
    2  = "Pharmaceuticals & Biotechnology"
    3  = "Auto Components"
    4  = "Minerals & Mining"
    5  = "Construction"
    6  = "Petroleum Products"
    7  = "Power"
    8  = "Finance"
    9  = "Insurance"
    10 = "Textiles & Apparels"
    11 = "Food Products"
    12 = "Telecom - Services"
    13 = "Agricultural Food & other Products"
    14 = "Agricultural, Commercial & Construction Vehicles"
    15 = "IT - Software"
    16 = "Power"
    17 = "Diversified"
    18 = "Transport Infrastructure"
    19 = "Electrical Equipment"
    20 = "Insurance"
    21 = "Power"
    22 = "Personal Products"
    23 = "Personal Products"
    24 = "Finance"
    25 = "IT - Software"
    26 = "Beverages"
    27 = "Realty"
    28 = "Transport Infrastructure"
    29 = "Diversified FMCG"
    30 = "Retailing"
    31 = "Construction"
    32 = "Telecom - Services"
    33 = "Finance"
    34 = "Pharmaceuticals & Biotechnology"
}

foreach ($row in $industries.Keys) {
    $ws.Range("C$row").Value = $industries[$row]
}
